$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New identifiers used throughout the report (replacing the old handback
# run's identifiers now that a fresh handoff has been generated).
# ---------------------------------------------------------------------------
$newMd1   = "320188b0-d72c-41c3-9bb7-f41b58c0e308.md"
$newMd2   = "ffffd236610a-41b7-4d3b-8d11-b4455e531a4e.md"
$newXlfZh = "320188b0-d72c-41c3-9bb7-f41b58c0e308.d86f4dfb38f8d79df51a423d2d5f3f6c9a795df5.zh-cn.xlf"
$newXlfDe = "320188b0-d72c-41c3-9bb7-f41b58c0e308.d86f4dfb38f8d79df51a423d2d5f3f6c9a795df5.de-de.xlf"

$statusText   = "Ready for handoff"
$overviewDate = "2016-03-25 07:58:58"
$zhHandoffDt  = "2016-03-25 07:58:53"
$deHandoffDt  = "2016-03-25 07:58:58"
$emptyDate    = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("D2").Value = $overviewDate

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
$wsOverview.Range("D3").Value = $overviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/320188b0-d72c-41c3-9bb7-f41b58c0e308.md", [Type]::Missing, [Type]::Missing, $newMd1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/ffffd236610a-41b7-4d3b-8d11-b4455e531a4e.md", [Type]::Missing, [Type]::Missing, $newMd2) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $zhHandoffDt
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("H2").Value = $emptyDate

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $zhHandoffDt
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()
$wsZh.Range("H3").Value = $emptyDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/320188b0-d72c-41c3-9bb7-f41b58c0e308.md", [Type]::Missing, [Type]::Missing, $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67537d25bf2ac946fba2a4af337f04ebb2e0031e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", [Type]::Missing, [Type]::Missing, $newXlfZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/ffffd236610a-41b7-4d3b-8d11-b4455e531a4e.md", [Type]::Missing, [Type]::Missing, $newMd2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67537d25bf2ac946fba2a4af337f04ebb2e0031e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", [Type]::Missing, [Type]::Missing, $newXlfZh) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $deHandoffDt
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("H2").Value = $emptyDate

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $deHandoffDt
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
$wsDe.Range("H3").Value = $emptyDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/320188b0-d72c-41c3-9bb7-f41b58c0e308.md", [Type]::Missing, [Type]::Missing, $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fe42363582669f540887a9fe8c4b18cfc5e283eb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", [Type]::Missing, [Type]::Missing, $newXlfDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/ffffd236610a-41b7-4d3b-8d11-b4455e531a4e.md", [Type]::Missing, [Type]::Missing, $newMd2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fe42363582669f540887a9fe8c4b18cfc5e283eb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", [Type]::Missing, [Type]::Missing, $newXlfDe) | Out-Null
